# Add March 2022 data for the executive board statistics sheet.
# Fills in Circulation / ILL Loans / ILL Borrows columns (B:D) for each
# library row (3-58) plus the Total row (59), which were previously blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(3,17223,2403,3259),
    @(4,7417,1010,1171),
    @(5,26949,2079,2802),
    @(6,370,203,7),
    @(7,16018,3367,2795),
    @(8,2162,425,440),
    @(9,2452,389,395),
    @(10,674,136,34),
    @(11,24,87,1),
    @(12,0,0,0),
    @(13,448,123,118),
    @(14,1032,635,516),
    @(15,1926,939,441),
    @(16,1362,649,330),
    @(17,690,258,50),
    @(18,6192,1043,1257),
    @(19,1209,330,242),
    @(20,8168,893,1552),
    @(21,154,170,8),
    @(22,7761,805,1358),
    @(23,413,244,57),
    @(24,6444,1327,1384),
    @(25,27176,3206,3420),
    @(26,1959,675,262),
    @(27,0,0,0),
    @(28,2144,510,479),
    @(29,1099,160,238),
    @(30,5676,1155,1003),
    @(31,211,124,90),
    @(32,729,695,105),
    @(33,6114,1248,1189),
    @(34,4492,1069,1254),
    @(35,2391,338,614),
    @(36,18780,2139,2359),
    @(37,3408,1091,606),
    @(38,6458,823,1131),
    @(39,288,305,61),
    @(40,677,98,188),
    @(41,1240,179,118),
    @(42,6049,226,160),
    @(43,218,82,70),
    @(44,373,45,18),
    @(45,1254,125,49),
    @(46,1221,357,173),
    @(47,4610,1304,950),
    @(48,11890,1316,2174),
    @(49,4880,1365,463),
    @(50,4068,467,742),
    @(51,9307,1276,1343),
    @(52,1608,198,339),
    @(53,5008,1115,1142),
    @(54,844,295,495),
    @(55,766,564,127),
    @(56,1470,416,619),
    @(57,4768,2346,1405),
    @(58,9134,657,415),
    @(59,245496,40481,40198)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
